$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 13, 14, 17, 18, 19
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F13").Value = 13548
$wsExhibit.Range("F14").Value = 181
$wsExhibit.Range("F17").Value = 5568
$wsExhibit.Range("F18").Value = 5591
$wsExhibit.Range("F19").Value = 61

# Sheet "全部类型" (All Types) - rows 35, 36, 40, 41, 42 (same events duplicated)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F35").Value = 13548
$wsAll.Range("F36").Value = 181
$wsAll.Range("F40").Value = 5568
$wsAll.Range("F41").Value = 5591
$wsAll.Range("F42").Value = 61
